# Update settlement recommendations:
# Bergi's recorded shared spend moves from 29.88 -> 30, and Adrien now shows
# a shared spend of 25 (was 0). Balances / Owes Matrix / Simplified
# Settlements are all formula-driven off these two "Shared Expenses" cells,
# so letting Excel recalculate takes care of the rest of the workbook.

$wb = $excel.ActiveWorkbook

$wsShared = $wb.Worksheets.Item("Shared Expenses")
$wsShared.Range("B2").Value = 30
$wsShared.Range("B4").Value = 25

# Mirror the cursor/selection state left behind in the author's session.
$wsDirect   = $wb.Worksheets.Item("Direct Expenses")
$wsBalances = $wb.Worksheets.Item("Balances")
$wsOwes     = $wb.Worksheets.Item("Owes Matrix")

$null = $wsDirect.Range("F8").Select()
$null = $wsBalances.Range("F3").Select()
$null = $wsOwes.Range("B17:C22").Select()

$null = $wsShared.Activate()
$null = $wsShared.Range("B4").Select()
